$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"0.002939958545068731"
$ws.Range("C2").Value = [double]"0.00293995854507047"
$ws.Range("D2").Value = [double]"0.002939958545078245"

$ws.Range("B3").Value = [double]"3.808932358547125E-07"
$ws.Range("C3").Value = [double]"3.808932362616549E-07"
$ws.Range("D3").Value = [double]"3.808932385690516E-07"

$ws.Range("B4").Value = [double]"1.448975640375648E-07"
$ws.Range("C4").Value = [double]"1.448975663858001E-07"
$ws.Range("D4").Value = [double]"1.448975633286896E-07"

$ws.Range("B5").Value = [double]"0.003667630002616969"
$ws.Range("C5").Value = [double]"0.003667630002608991"
$ws.Range("D5").Value = [double]"0.00366763000261794"

$ws.Range("B6").Value = [double]"1.366435165840219E-05"
$ws.Range("C6").Value = [double]"1.366435169617108E-05"
$ws.Range("D6").Value = [double]"1.366435169313481E-05"

$ws.Range("B7").Value = [double]"2.368757855338019E-05"
$ws.Range("C7").Value = [double]"2.209255172651603E-05"
$ws.Range("D7").Value = [double]"2.368757850949929E-05"

$ws.Range("B8").Value = [double]"0.0005482768665219502"
$ws.Range("C8").Value = [double]"0.000434670253791872"
$ws.Range("D8").Value = [double]"0.0006426513249862472"

$ws.Range("B9").Value = [double]"0.001551349602901047"
$ws.Range("C9").Value = [double]"0.00155134960290326"
$ws.Range("D9").Value = [double]"0.001551349602902931"

$ws.Range("B10").Value = [double]"0.000834135881859916"
$ws.Range("C10").Value = [double]"0.0005011916252527399"
$ws.Range("D10").Value = [double]"2.441486563903493E-05"

$ws.Range("B11").Value = [double]"6.852026537000017E-05"
$ws.Range("D11").Value = [double]"6.852027282981138E-05"
